$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119, shifting existing rows 119:136 down to 120:137.
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new weekly record.
$ws.Cells.Item(119, 1).Value = 7
$ws.Cells.Item(119, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(119, 3).Value = "Ñuble"
$ws.Cells.Item(119, 4).Value = 44504
$ws.Cells.Item(119, 5).Value = 16
$ws.Cells.Item(119, 6).Value = 100112032
$ws.Cells.Item(119, 7).Value = "Zapallo italiano"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 100
$ws.Cells.Item(119, 11).Value = 14000
$ws.Cells.Item(119, 12).Value = 15000
$ws.Cells.Item(119, 13).Value = 14500
$ws.Cells.Item(119, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(119, 15).Value = "Región del Maule"
$ws.Cells.Item(119, 16).Value = 242
$ws.Cells.Item(119, 17).Value = 60
$ws.Cells.Item(119, 18).Value = "Hortaliza"
